# The ClinvarID column (I) on both sheets currently stores pipe-delimited
# RCV accession lists. Replace them with a single numeric ClinVar "VCV"-style
# identifier, keeping the cell content as TEXT (not a number) since that is
# how every other numeric-looking identifier column in this workbook
# (ACMG_version, OMIM_disorder, ...) is stored.
#
# Writing a pure-digit string straight into .Value would make Excel infer a
# Number type. To avoid that we briefly mark the cell as Text before the
# write (so the literal is kept verbatim) and then reset the cell style back
# to Normal so no stray number-format/style is left behind on the cell.
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$wsPR = $wb.Worksheets.Item("PR results")
$wsRR = $wb.Worksheets.Item("RR results")

Set-TextValue $wsPR.Range("I2") "89363"

Set-TextValue $wsRR.Range("I2") "866323"
Set-TextValue $wsRR.Range("I3") "43492"
Set-TextValue $wsRR.Range("I4") "7105"
Set-TextValue $wsRR.Range("I5") "1677653"
Set-TextValue $wsRR.Range("I6") "92734"
Set-TextValue $wsRR.Range("I7") "374315"
